$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.457.14'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '2.249.78'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  -0.38%  '
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.03'
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("E6").Value = '  +0.83%  '
$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.50'
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = '  +8.31%  '
$ws.Range("E8").Value = '  -0.45%  '
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.639'
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = '  +0.69%  '
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.20'
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = '  +5.65%  '
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0955'
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = '  -0.19%  '
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.26'
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = '  +0.80%  '
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").Value = '2.590.51'
$ws.Range("E14").Value = '  -0.85%  '
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.91'
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = '  +1.08%  '
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.862'
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").Value = '2.269.12'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '42.393.14'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '0.0₃0982'
$ws.Range("E19").Value = '  -0.20%  '
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.16'
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = '  -0.80%  '
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.57'
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = '  -1.22%  '
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.58'
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("E25").Value = '  -4.30%  '
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.21'
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = '  -1.41%  '
$ws.Range("E27").Value = '  -4.27%  '
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.99'
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = '  +11.20%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.14'
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = '  +1.22%  '
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.60'
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = '  -0.88%  '
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0854'
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = '  +8.46%  '
$ws.Range("E33").Value = '  -3.38%  '
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.90'
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("E35").Value = '  +1.25%  '
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.48'
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = '  -3.63%  '
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.73'
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = '  +0.89%  '
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0299'
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = '  -1.98%  '
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.05'
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = '  -1.45%  '
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.23'
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = '  -2.23%  '
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.94'
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = '  +0.04%  '
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '116.91'
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = '  +23.51%  '
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.205'
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = '  -0.71%  '
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.39'
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = '  -1.14%  '
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.77'
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = '  -4.56%  '
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("E47").Value = '  -0.51%  '
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.13'
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.17'
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("B50").Value = 'FTXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.31'
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = '  -11.39%  '
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.13'
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = '  -0.76%  '
